$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: add the new row 8 using row 7's CURRENT formatting (s=4/s=5) ---
# Copy row 7's formats down to row 8 before we change row 7's own formatting.
# (Only B:E -- row 7, like rows 5/6, has no populated/styled A cell.)
$ws.Range("B7:E7").Copy()
$ws.Range("B8:E8").PasteSpecial(-4122)

# --- Step 2: re-style row 7 so it matches row 3's "section separator" look (s=6/s=7) ---
$ws.Range("A3:E3").Copy()
$ws.Range("A7:E7").PasteSpecial(-4122)

# --- Step 3: populate row 8 values (the new Ksat dialogue line) ---
# NOTE: the source data uses a literal backslash-n (two characters), not a
# real line break, so these are single-quoted (no backtick escape expansion).
$ws.Range("B8").Value = 188
$ws.Range("C8").Value = ' ...[K]You two... I commend you\nfor graduating.'
$ws.Range("D8").Value = ' ...[K]Вы двое... Поздравляю вас\nс выпуском.'
$ws.Range("E8").Value = ' ...[K]Âú äâïå... Ðïèäñàâìÿý âàò\nò âúðôòëïí.'
$ws.Rows.Item(8).RowHeight = 21.6

# --- Step 4: update the saved selection/view state ---
$excel.Goto($ws.Range("D12"), $true)
